# Update quizvragen via Admin
#
# DC sheet: remove the "q1" (Ohmse wet) row and the trailing blank "new
# question" row, leaving only the "q3" (Vermogen) row, which becomes row 2.
#
# Wiskunde 3 sheet: a new question row is inserted at the top (row 2),
# pushing every existing question down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "DC"
# ---------------------------------------------------------------------
$dc = $wb.Worksheets.Item("DC")

# Delete the last (blank test) row first so row indices for the row we
# delete next don't shift.
$dc.Rows.Item(4).Delete()

# Delete the old row 2 ("q1" Ohmse wet); the old row 3 ("q3" Vermogen)
# shifts up to become the new row 2.
$dc.Rows.Item(2).Delete()

# ---------------------------------------------------------------------
# Sheet "Wiskunde 3"
# ---------------------------------------------------------------------
$wk = $wb.Worksheets.Item("Wiskunde 3")

# Insert a new blank row at position 2; existing rows 2-20 shift down to 3-21.
$wk.Rows.Item(2).Insert()

# Insert copies formatting from the row above (the bold header); reset the
# new row back to the plain style used by the other data rows.
$wk.Range("A2:L2").Style = "Normal"

$wk.Range("A2").Value = "q1"
$wk.Range("B2").Value = "mc"
$wk.Range("C2").Value = "Goniometrie"
$wk.Range("D2").Value = "Wat is de juiste formule voor de stroom I?"
$wk.Range("E2").Value = "['I = U/R', ' U = I*R', ' R = U/I']"
$wk.Range("F2").Value = 0
$wk.Range("G2").Value = "sin(α) = overstaande / schuine"
# H2 (image_path) stays blank, like the freshly inserted cell.
$wk.Range("I2").Value = "sin(α)=o/h"
$wk.Range("J2").Value = '["sinus","basisformule"]'
$wk.Range("K2").Value = 1
$wk.Range("L2").Value = "https://raw.githubusercontent.com/onomatorHanze/didactic-octo-spork/main/data/images/Wiskunde_3_q0_1763116823.jpg"
